$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and G (Hora) hold numeric-looking text (e.g. "247.52", "6").
# Force each target cell to Text format before assigning the new value so Excel
# keeps it as a string (preserving exact digits/precision) instead of silently
# converting it to a number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G51").NumberFormat = "@"

$ws.Range("D2").Value = "247.52"
$ws.Range("E2").Value = "1BNBBNB"
$ws.Range("G2").Value = "6"
$ws.Range("D3").Value = "22.68"
$ws.Range("G3").Value = "6"
$ws.Range("D4").Value = "5.625"
$ws.Range("G4").Value = "6"
$ws.Range("D5").Value = "0.05620"
$ws.Range("G5").Value = "6"
$ws.Range("D6").Value = "3.403"
$ws.Range("G6").Value = "6"
$ws.Range("D7").Value = "6.467"
$ws.Range("G7").Value = "6"
$ws.Range("G8").Value = "6"
$ws.Range("G9").Value = "6"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1424"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("G10").Value = "6"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07477"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("G11").Value = "6"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03193"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").Value = "6"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.02980"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("G13").Value = "6"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09263"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("G14").Value = "6"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001664"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("G15").Value = "6"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.252"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("G16").Value = "6"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04692"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("G17").Value = "6"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005754"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("G18").Value = "6"
$ws.Range("D19").Value = "0.006267"
$ws.Range("G19").Value = "6"
$ws.Range("D20").Value = "0.001057"
$ws.Range("G20").Value = "6"
$ws.Range("D21").Value = "0.003811"
$ws.Range("G21").Value = "6"
$ws.Range("D22").Value = "0.0001501"
$ws.Range("G22").Value = "6"
$ws.Range("D23").Value = "0.0004603"
$ws.Range("G23").Value = "6"
$ws.Range("D24").Value = "3.980"
$ws.Range("G24").Value = "6"
$ws.Range("D25").Value = "2.122"
$ws.Range("G25").Value = "6"
$ws.Range("G26").Value = "6"
$ws.Range("E27").Value = "26ProBitTokenPROBBestin24h"
$ws.Range("G27").Value = "6"
$ws.Range("G28").Value = "6"
$ws.Range("G29").Value = "6"
$ws.Range("G30").Value = "6"
$ws.Range("G31").Value = "6"
$ws.Range("G32").Value = "6"
$ws.Range("G33").Value = "6"
$ws.Range("G34").Value = "6"
$ws.Range("G35").Value = "6"
$ws.Range("G36").Value = "6"
$ws.Range("G37").Value = "6"
$ws.Range("G38").Value = "6"
$ws.Range("G39").Value = "6"
$ws.Range("D40").Value = "0.04192"
$ws.Range("G40").Value = "6"
$ws.Range("D41").Value = "0.007145"
$ws.Range("G41").Value = "6"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1045"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("G42").Value = "6"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.003303"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("G43").Value = "6"
$ws.Range("D44").Value = "0.009798"
$ws.Range("G44").Value = "6"
$ws.Range("D45").Value = "0.00005684"
$ws.Range("G45").Value = "6"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("G46").Value = "6"
$ws.Range("D47").Value = "0.6804"
$ws.Range("G47").Value = "6"
$ws.Range("D48").Value = "0.02900"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("G48").Value = "6"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("G49").Value = "6"
$ws.Range("D50").Value = "0.01011"
$ws.Range("G50").Value = "6"
$ws.Range("G51").Value = "6"
